$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 78 (Z16_B01_P01_Ib01_I01 / "16.1 Anzahl an Straftaten") entirely.
# All subsequent rows shift up by one, and the sheet dimension shrinks from
# A1:M84 to A1:M83.
$ws.Rows.Item(78).Delete()
